$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.118.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.991.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4975"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.74%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4198"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09206"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.096"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.013.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.991"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.452"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.013"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06761"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.011"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.977"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.116.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.289"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.250.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.80%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.282"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.261"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.048"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09851"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.536"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.823"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.744"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02426"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.320"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.072"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06402"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6484"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.79%  "

$ws.Range("E43").Value = "  -3.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.010"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.34%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6214"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.345"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.189"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.491"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000338"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.24%  "

$ws.Range("E51").Value = "  -0.22%  "
